# Updated loading_percent results for the 380 kV case (rows 2-25, data rows 0-23)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","D","E","F","G","H","I","K","L","N","O")

$newValues = @{
    2 = @(13.24368449333586,4.337361940648294,16.84989591380234,20.90432518335481,22.95729994256521,12.8213740005175,24.15813738338176,10.27389213124116,8.566822695259878,18.51634187932808,18.73949653224093)
    3 = @(13.09908436810551,4.252353523591249,16.90999013700898,20.90163755142451,22.94292176482539,12.85609689819221,24.26171910333506,9.992388720086332,8.523149269659488,18.56724572392944,18.78500994135695)
    4 = @(13.01176250836271,4.198551107563008,16.94938549518691,20.9055191971832,22.94176770243483,12.87933789514048,24.32920374005498,9.81378304373689,8.497685128956274,18.6001137216249,18.81689362474793)
    5 = @(12.97658393914115,4.176237352236346,16.96606811446764,20.90849330204366,22.94322799543792,12.88929189257574,24.35768225928737,9.739617831980595,8.487656415489422,18.61391427322798,18.83087534297428)
    6 = @(12.9707681665252,4.172509146373072,16.9688762438348,20.90907126008923,22.94358706693141,12.89097392579934,24.36247019500724,9.727221346276306,8.486012414461202,18.61623043047996,18.83325667574419)
    7 = @(13.01128638530501,4.198251730124563,16.94960793628448,20.90555366864527,22.94177957985786,12.87947018210655,24.32958385008392,9.812788329014841,8.497548458050549,18.60029819312589,18.81707818527511)
    8 = @(13.19354559312829,4.308393439166779,16.87009853694783,20.90225139211925,22.95075062659974,12.83294778789675,24.19304665418053,10.17806396785093,8.551488693006835,18.53355939450267,18.75437134221456)
    9 = @(13.56087524913835,4.511004335914559,16.73396230285814,20.93957223615115,23.02911343592013,12.75695722577475,23.95607541693242,10.84599127654855,8.667605346471174,18.41543519386444,18.6627110635667)
    10 = @(13.83451619385099,4.65092969959322,16.64595628592924,20.99350266056049,23.1234493025956,12.71041364603971,23.80066802454267,11.30403181609563,8.758679313211221,18.33635644764661,18.6145267774034)
    11 = @(13.95933517427322,4.712498016448287,16.60851906062887,21.02373261370031,23.17424689243916,12.69125525512587,23.73401493414901,11.50473418890414,8.80124015327614,18.30204108056358,18.59677846294728)
    12 = @(14.00660943796252,4.735501330120637,16.59471534308652,21.03599227551651,23.19460472459778,12.68429006892462,23.70935555392321,11.57959129146671,8.817508672190813,18.28928413670351,18.59065804466384)
    13 = @(13.99642836208732,4.730561173546151,16.59767164186959,21.03331593498239,23.19017061263372,12.6857772630317,23.71464057523398,11.56352101089756,8.813998387758932,18.29202102183241,18.59194947065621)
    14 = @(13.96322452585241,4.714396814333885,16.607375948997,21.02472498844271,23.17589933345065,12.69067641952183,23.73197455738424,11.51091589458224,8.802575598384935,18.30098680337725,18.59626289459759)
    15 = @(13.94288606822078,4.704454812868077,16.6133686708252,21.01956833913748,23.16730348376302,12.69371501887069,23.74266772863328,11.47854349948939,8.795598230393059,18.30650951269763,18.59898320921646)
    16 = @(13.82636230039674,4.646863116100828,16.6484551132478,20.99164112393271,23.12028723747055,12.71170622919351,23.80510525340041,11.29075750509793,8.755919686135574,18.33863232507785,18.6157706669114)
    17 = @(13.7549362947969,4.610990128990146,16.67064432826932,20.97596312919307,23.09345572576078,12.72325922247052,23.84444359036013,11.17356253694635,8.731860098941615,18.35876264733074,18.62713807944309)
    18 = @(13.7138881882589,4.590161374437965,16.68365146456069,20.96748224515095,23.07876583121632,12.73009381584299,23.86745049564172,11.10543611209122,8.718128773119449,18.37049718163123,18.63406887599618)
    19 = @(13.69999713428428,4.583075893424326,16.68809746429293,20.96470313467017,23.07392000384242,12.73244045851703,23.87530561482442,11.08224756296332,8.713498308727919,18.37449713132024,18.63648292251731)
    20 = @(13.76253646741779,4.614829195754901,16.66825694881236,20.97757658086126,23.09623517989695,12.72200976191928,23.84021657790452,11.18611287608119,8.734410275974616,18.3566035901898,18.62588736455026)
    21 = @(13.97297740003204,4.719153214039629,16.60451544116381,21.02722637424105,23.18006080965002,12.68922955677678,23.72686739131603,11.52639867278709,8.805926721216572,18.29834689607983,18.59497963592964)
    22 = @(14.11054116527491,4.785515899113601,16.56503027016018,21.06440577946436,23.24137813614813,12.66949426351676,23.65617164619306,11.74210717616096,8.853545560211913,18.26165695494869,18.57827989906836)
    23 = @(14.03713142828785,4.750266926363346,16.58590553414052,21.04413215449667,23.20805858484335,12.67987287229665,23.69359380495335,11.62760440333908,8.828053766795131,18.281112698814,18.5868724124913)
    24 = @(13.75910037750631,4.613094190996812,16.66933550350226,20.97684547983967,23.09497629555607,12.72257404302494,23.84212639297209,11.18044119895266,8.733257025970781,18.35757919728134,18.62645158079722)
    25 = @(13.46067635182941,4.457708416153247,16.76867782775995,20.92480274107642,23.00143005333075,12.77588351440136,24.01689457252221,10.67082790829891,8.635140051237727,18.44603276441675,18.6627110635667)
}

foreach ($row in $newValues.Keys) {
    $rowValues = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $rowValues[$i]
    }
}
